# SwaadSutra_Consolidated_2026-01-13.xlsx update
# A new order (#12 "Vermicelli Kheer x1" from Swapnil (Phantom)) was placed
# after order #11. The "All Orders" sheet keeps newest orders on top, so a
# new row is inserted at row 2 (pushing every existing order down by one
# row) and populated with the new order's details. The "Daily Summary"
# sheet is updated to reflect the new totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("All Orders")

# Insert a new blank row at the top of the data (row 2), shifting all
# existing order rows down by one.
$ws.Rows.Item(2).Insert()

function Set-TextCell {
    param($Cell, [string]$Text)
    # Force the cell to be stored as text (even when the text looks like a
    # number or a date, e.g. "420" or "2026-01-15"), matching the source
    # data's convention of keeping every non Order-ID/Total column as text.
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    # Drop back to the default "Normal" style so no explicit cell style is
    # left behind (the workbook only uses the default style).
    $Cell.Style = "Normal"
}

# New order field values
$orderId       = 12
$orderDate     = "2026-01-13 22:43"
$customer      = "Swapnil (Phantom)"
$flatNo        = "420"
$phone         = ""
$items         = "Vermicelli Kheer x1"
$total         = 50
$status        = "NEW"
$payment       = "PENDING"
$collectionDate = "2026-01-15"
$collectionTime = "16:42"
$notes         = "No vermicelli in kheer please."
$cancelReason  = ""
$feedback      = ""

$ws.Cells.Item(2,1).Value = $orderId
Set-TextCell $ws.Cells.Item(2,2)  $orderDate
Set-TextCell $ws.Cells.Item(2,3)  $customer
Set-TextCell $ws.Cells.Item(2,4)  $flatNo
Set-TextCell $ws.Cells.Item(2,5)  $phone
Set-TextCell $ws.Cells.Item(2,6)  $items
$ws.Cells.Item(2,7).Value = $total
Set-TextCell $ws.Cells.Item(2,8)  $status
Set-TextCell $ws.Cells.Item(2,9)  $payment
Set-TextCell $ws.Cells.Item(2,10) $collectionDate
Set-TextCell $ws.Cells.Item(2,11) $collectionTime
Set-TextCell $ws.Cells.Item(2,12) $notes
Set-TextCell $ws.Cells.Item(2,13) $cancelReason
Set-TextCell $ws.Cells.Item(2,14) $feedback

# Update the Daily Summary sheet: total orders, revenue, and pending
# payment amounts all grow by the new order's 50 rupee total.
$summary = $wb.Worksheets.Item("Daily Summary")
$summary.Cells.Item(2,2).Value = 12   # Total Orders
$summary.Cells.Item(2,5).Value = 325  # Revenue
$summary.Cells.Item(2,7).Value = 325  # Pending
